# Updates the Optostim Receiver BOM worksheet with Digi-Key pricing / packaging
# details for line items 3-11 (rows 6-14), which previously only had the
# Manufacturer Part Number / Description / lookup columns populated.
#
# Columns:
#   H = Requested Quantity 1   (number)
#   I = Pack Quantity 1        (number)
#   J = Pack Type 1            (text)
#   K = Digi-Key Part Number 1 (text)
#   L = Unit Price 1           (text, keeps trailing zeros e.g. "0.57000")
#   M = Extended Price 1       (text, currency-formatted e.g. "$0.57")
#   N = Minimum Order Quantity 1 (number)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Plain($range, $value) {
    $range.Value = $value
}

function Set-TextValue($range, $value) {
    # Force the value to be stored as text even though it looks numeric
    # (e.g. "0.57000", "$0.57", "1,558,028") so Excel does not silently
    # convert it to a number and lose formatting / precision.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-Plain $ws.Range("H6") 1
Set-Plain $ws.Range("I6") 1
Set-TextValue $ws.Range("J6") "Tube"
Set-TextValue $ws.Range("K6") "2223-DS01C-254-L-03BE-ND"
Set-TextValue $ws.Range("L6") "0.57000"
Set-TextValue $ws.Range("M6") "`$0.57"
Set-Plain $ws.Range("N6") 1

Set-Plain $ws.Range("H7") 1
Set-Plain $ws.Range("I7") 1
Set-TextValue $ws.Range("J7") "Bulk"
Set-TextValue $ws.Range("K7") "A123828-ND"
Set-TextValue $ws.Range("L7") "2.41000"
Set-TextValue $ws.Range("M7") "`$2.41"
Set-Plain $ws.Range("N7") 1

Set-Plain $ws.Range("H8") 3
Set-Plain $ws.Range("I8") 3
Set-TextValue $ws.Range("J8") "Cut Tape (CT)"
Set-TextValue $ws.Range("K8") "CR0805-FX-4701ELFCT-ND"
Set-TextValue $ws.Range("L8") "0.10000"
Set-TextValue $ws.Range("M8") "`$0.30"
Set-Plain $ws.Range("N8") 1

Set-Plain $ws.Range("H9") 1
Set-Plain $ws.Range("I9") 1
Set-TextValue $ws.Range("J9") "Tray"
Set-TextValue $ws.Range("K9") "360-3252-ND"
Set-TextValue $ws.Range("L9") "4.97000"
Set-TextValue $ws.Range("M9") "`$4.97"
Set-Plain $ws.Range("N9") 1

Set-Plain $ws.Range("H10") 1
Set-Plain $ws.Range("I10") 1
Set-TextValue $ws.Range("J10") "Cut Tape (CT)"
Set-TextValue $ws.Range("K10") "NCP1117ST50T3GOSCT-ND"
Set-TextValue $ws.Range("L10") "0.72000"
Set-TextValue $ws.Range("M10") "`$0.72"
Set-Plain $ws.Range("N10") 1

Set-Plain $ws.Range("H11") 2
Set-Plain $ws.Range("I11") 2
Set-TextValue $ws.Range("J11") "Cut Tape (CT)"
Set-TextValue $ws.Range("K11") "1276-1052-1-ND"
Set-TextValue $ws.Range("L11") "0.11000"
Set-TextValue $ws.Range("M11") "`$0.22"
Set-Plain $ws.Range("N11") 1

Set-Plain $ws.Range("H12") 1
Set-Plain $ws.Range("I12") 1
Set-TextValue $ws.Range("J12") "Tray"
Set-TextValue $ws.Range("K12") "CP-102A-ND"
Set-TextValue $ws.Range("L12") "0.77000"
Set-TextValue $ws.Range("M12") "`$0.77"
Set-Plain $ws.Range("N12") 1

Set-Plain $ws.Range("H13") 1
Set-Plain $ws.Range("I13") 1
Set-TextValue $ws.Range("J13") "Cut Tape (CT)"
Set-TextValue $ws.Range("K13") "CR0805-FX-1001ELFCT-ND"
Set-TextValue $ws.Range("L13") "0.10000"
Set-TextValue $ws.Range("M13") "`$0.10"
Set-Plain $ws.Range("N13") 1

Set-Plain $ws.Range("H14") 1
Set-Plain $ws.Range("I14") 1
Set-TextValue $ws.Range("J14") "Tray"
Set-TextValue $ws.Range("K14") "WM5514-ND"
Set-TextValue $ws.Range("L14") "2.65000"
Set-TextValue $ws.Range("M14") "`$2.65"
Set-Plain $ws.Range("N14") 1

# Row 13 (CR0805-FX-1001ELF): Availability quantity increased.
Set-TextValue $ws.Range("E13") "1,558,028"
